$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "last checked" timestamp batch (rows 968-981)
for ($r = 968; $r -le 981; $r++) {
  $ws.Cells.Item($r, 4).Value = 44233.13005762731
}

# Row 982
$ws.Range("A982").Value = "Odoo"
$b982 = $ws.Range("B982")
$b982.Value = "https://www.dataintelligence-group.com/"
$b982.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b982, "https://www.dataintelligence-group.com/")
$b982.Style = "Hyperlink"
$ws.Range("C982").Value = "Disponible"
$d982 = $ws.Range("D982")
$d982.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d982.Value = 44233.1512524904

# Row 983
$ws.Range("A983").Value = "Blackbox"
$b983 = $ws.Range("B983")
$b983.Value = "https://serviciodashboard.azurewebsites.net/"
$b983.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b983, "https://serviciodashboard.azurewebsites.net/")
$b983.Style = "Hyperlink"
$ws.Range("C983").Value = "Disponible"
$d983 = $ws.Range("D983")
$d983.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d983.Value = 44233.1512524904

# Row 984
$ws.Range("A984").Value = "PowerBI"
$b984 = $ws.Range("B984")
$b984.Value = "https://powerbi.microsoft.com/es-es/"
$b984.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b984, "https://powerbi.microsoft.com/es-es/")
$b984.Style = "Hyperlink"
$ws.Range("C984").Value = "Disponible"
$d984 = $ws.Range("D984")
$d984.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d984.Value = 44233.1512524904

# Row 985
$ws.Range("A985").Value = "Dropbox"
$b985 = $ws.Range("B985")
$b985.Value = "https://www.dropbox.com/"
$b985.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b985, "https://www.dropbox.com/")
$b985.Style = "Hyperlink"
$ws.Range("C985").Value = "Disponible"
$d985 = $ws.Range("D985")
$d985.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d985.Value = 44233.1512524904

# Row 986
$ws.Range("A986").Value = "Odoo"
$b986 = $ws.Range("B986")
$b986.Value = "https://dataintelligence.store/"
$b986.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b986, "https://dataintelligence.store/")
$b986.Style = "Hyperlink"
$ws.Range("C986").Value = "Disponible"
$d986 = $ws.Range("D986")
$d986.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d986.Value = 44233.1512524904

# Row 987
$ws.Range("A987").Value = "GEE"
$b987 = $ws.Range("B987")
$b987.Value = "https://app-data-i.users.earthengine.app/"
$b987.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b987, "https://app-data-i.users.earthengine.app/")
$b987.Style = "Hyperlink"
$ws.Range("C987").Value = "Disponible"
$d987 = $ws.Range("D987")
$d987.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d987.Value = 44233.1512524904

# Row 988
$ws.Range("A988").Value = "UtilidadesOdoo"
$b988 = $ws.Range("B988")
$b988.Value = "https://odooutil.azurewebsites.net/"
$b988.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b988, "https://odooutil.azurewebsites.net/")
$b988.Style = "Hyperlink"
$ws.Range("C988").Value = "Disponible"
$d988 = $ws.Range("D988")
$d988.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d988.Value = 44233.1512524904

# Row 989
$ws.Range("A989").Value = "Filtros Dashboard"
$b989 = $ws.Range("B989")
$b989.Value = "https://filtradordashboard.azurewebsites.net/"
$b989.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b989, "https://filtradordashboard.azurewebsites.net/")
$b989.Style = "Hyperlink"
$ws.Range("C989").Value = "Disponible"
$d989 = $ws.Range("D989")
$d989.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d989.Value = 44233.1512524904

# Row 990
$ws.Range("A990").Value = "MapStore"
$b990 = $ws.Range("B990")
$b990.Value = "https://ide.dataintelligence-group.com/mapstore/#/"
$b990.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b990, "https://ide.dataintelligence-group.com/mapstore/", "/")
$b990.Style = "Hyperlink"
$ws.Range("C990").Value = "Disponible"
$d990 = $ws.Range("D990")
$d990.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d990.Value = 44233.1512524904

# Row 991
$ws.Range("A991").Value = "GeoServer"
$b991 = $ws.Range("B991")
$b991.Value = "https://ide.dataintelligence-group.com/geoserver/web/?0"
$b991.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b991, "https://ide.dataintelligence-group.com/geoserver/web/?0")
$b991.Style = "Hyperlink"
$ws.Range("C991").Value = "Disponible"
$d991 = $ws.Range("D991")
$d991.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d991.Value = 44233.1512524904

# Row 992
$ws.Range("A992").Value = "Tomcat"
$b992 = $ws.Range("B992")
$b992.Value = "https://ide.dataintelligence-group.com/"
$b992.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b992, "https://ide.dataintelligence-group.com/")
$b992.Style = "Hyperlink"
$ws.Range("C992").Value = "Disponible"
$d992 = $ws.Range("D992")
$d992.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d992.Value = 44233.1512524904

# Row 993
$ws.Range("A993").Value = "Shiny"
$b993 = $ws.Range("B993")
$b993.Value = "https://rpubs.com/dataintelligence/"
$b993.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b993, "https://rpubs.com/dataintelligence/")
$b993.Style = "Hyperlink"
$ws.Range("C993").Value = "Disponible"
$d993 = $ws.Range("D993")
$d993.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d993.Value = 44233.1512524904

# Row 994
$ws.Range("A994").Value = "Github"
$b994 = $ws.Range("B994")
$b994.Value = "https://github.com/Sud-Austral/"
$b994.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b994, "https://github.com/Sud-Austral/")
$b994.Style = "Hyperlink"
$ws.Range("C994").Value = "Disponible"
$d994 = $ws.Range("D994")
$d994.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d994.Value = 44233.1512524904

# Row 995
$ws.Range("A995").Value = "EZ Exporter"
$b995 = $ws.Range("B995")
$b995.Value = "https://ezexporter.highviewapps.com/exports/export-profile/"
$b995.Style = "Hyperlink"
$null = $ws.Hyperlinks.Add($b995, "https://ezexporter.highviewapps.com/exports/export-profile/")
$b995.Style = "Hyperlink"
$ws.Range("C995").Value = "Disponible"
$d995 = $ws.Range("D995")
$d995.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d995.Value = 44233.1512524904

